$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H107").Value = 531.4
$ws.Range("I107").Value = 557.0833
$ws.Range("K107").Value = 557.0833
$ws.Range("M107").Value = 1362.9167
$ws.Range("H137").Value = 1887.7
$ws.Range("I137").Value = 1485.6428
$ws.Range("K137").Value = 4456.928400000001
$ws.Range("M137").Value = -1906.928400000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 387
$ws.Range("I94").Value = 376.58334
$ws.Range("K94").Value = 376.58334
$ws.Range("M94").Value = 74.41665999999998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 5015.294
$ws.Range("I22").Value = 10239.6
$ws.Range("K22").Value = 10239.6
$ws.Range("M22").Value = -9889.6
$ws.Range("H31").Value = 2464.423
$ws.Range("I31").Value = 1065.8182
$ws.Range("K31").Value = 1065.8182
$ws.Range("M31").Value = -770.8181999999999
$ws.Range("H34").Value = 2464.423
$ws.Range("I34").Value = 1065.8182
$ws.Range("K34").Value = 1065.8182
$ws.Range("M34").Value = -863.8181999999999
$ws.Range("H62").Value = 9814.929
$ws.Range("I62").Value = 6127.273
$ws.Range("J62").Value = 23336.334
$ws.Range("K62").Value = 6127.273
$ws.Range("L62").Value = 23336.334
$ws.Range("M62").Value = -5503.273
$ws.Range("N62").Value = -24584.334
$ws.Range("H65").Value = 9814.929
$ws.Range("I65").Value = 6127.273
$ws.Range("J65").Value = 23336.334
$ws.Range("K65").Value = 30636.365
$ws.Range("L65").Value = 116681.67
$ws.Range("M65").Value = -27516.365
$ws.Range("N65").Value = -122921.67
$ws.Range("H99").Value = 2394.75
$ws.Range("I99").Value = 2299.5
$ws.Range("K99").Value = 2299.5
$ws.Range("M99").Value = -801.5
$ws.Range("H126").Value = 2394.75
$ws.Range("I126").Value = 2299.5
$ws.Range("K126").Value = 6898.5
$ws.Range("M126").Value = -4428.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 48.333332
$ws.Range("I15").Value = 48.333332
$ws.Range("K15").Value = 144.999996
$ws.Range("M15").Value = -4.99999600000001
$ws.Range("H40").Value = 81.69231000000001
$ws.Range("I40").Value = 57.205883
$ws.Range("J40").Value = 127.94444
$ws.Range("K40").Value = 228.823532
$ws.Range("L40").Value = 511.77776
$ws.Range("M40").Value = -159.823532
$ws.Range("N40").Value = -649.7777599999999
$ws.Range("H59").Value = 2002.5
$ws.Range("I59").Value = 505
$ws.Range("J59").Value = 3500
$ws.Range("K59").Value = 1515
$ws.Range("L59").Value = 10500
$ws.Range("M59").Value = -975
$ws.Range("N59").Value = -11580
$ws.Range("H121").Value = 1401391.2
$ws.Range("I121").Value = 91800.63
$ws.Range("J121").Value = 5002765.5
$ws.Range("K121").Value = 275401.89
$ws.Range("L121").Value = 15008296.5
$ws.Range("M121").Value = -274091.89
$ws.Range("N121").Value = -15010916.5
$ws.Range("H132").Value = 996
$ws.Range("I132").Value = 994
$ws.Range("J132").Value = 1000
$ws.Range("K132").Value = 8946
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -6416
$ws.Range("N132").Value = -14060
$ws.Range("H133").Value = 7688
$ws.Range("I133").Value = 7688
$ws.Range("K133").Value = 23064
$ws.Range("M133").Value = -18004

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 112.27273
$ws.Range("I2").Value = 103.6
$ws.Range("K2").Value = 103.6
$ws.Range("M2").Value = 9.400000000000006
$ws.Range("H3").Value = 6011996.5
$ws.Range("I3").Value = 3343333
$ws.Range("J3").Value = 7155709.5
$ws.Range("K3").Value = 3343333
$ws.Range("L3").Value = 7155709.5
$ws.Range("M3").Value = -3343217
$ws.Range("N3").Value = -7155941.5
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 4221
$ws.Range("J7").Value = 4221
$ws.Range("L7").Value = 4221
$ws.Range("N7").Value = -4445
$ws.Range("H8").Value = 4221
$ws.Range("J8").Value = 4221
$ws.Range("L8").Value = 4221
$ws.Range("N8").Value = -4499
$ws.Range("H9").Value = 3594
$ws.Range("I9").Value = 2101.818
$ws.Range("J9").Value = 20008
$ws.Range("K9").Value = 2101.818
$ws.Range("L9").Value = 20008
$ws.Range("M9").Value = -1931.818
$ws.Range("N9").Value = -20348
$ws.Range("H10").Value = 1000501
$ws.Range("I10").Value = 1500000
$ws.Range("J10").Value = 1503
$ws.Range("K10").Value = 1500000
$ws.Range("L10").Value = 1503
$ws.Range("M10").Value = -1499831
$ws.Range("N10").Value = -1841
$ws.Range("H11").Value = 3691412
$ws.Range("I11").Value = 4217407.5
$ws.Range("J11").Value = 9443.5
$ws.Range("K11").Value = 4217407.5
$ws.Range("L11").Value = 9443.5
$ws.Range("M11").Value = -4217268.5
$ws.Range("N11").Value = -9721.5
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("H13").Value = 3214.2856
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 6966.6665
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 6966.6665
$ws.Range("M13").Value = -261
$ws.Range("N13").Value = -7244.6665
$ws.Range("H14").Value = 15683333
$ws.Range("I14").Value = 18800000
$ws.Range("J14").Value = 99999
$ws.Range("K14").Value = 18800000
$ws.Range("L14").Value = 99999
$ws.Range("M14").Value = -18799832
$ws.Range("N14").Value = -100335
$ws.Range("H17").Value = 11029.5
$ws.Range("J17").Value = 11029.5
$ws.Range("L17").Value = 11029.5
$ws.Range("N17").Value = -11365.5
$ws.Range("H18").Value = 349999
$ws.Range("J18").Value = 24998.5
$ws.Range("L18").Value = 24998.5
$ws.Range("N18").Value = -25584.5
$ws.Range("H19").Value = 8300.333000000001
$ws.Range("I19").Value = 2450
$ws.Range("J19").Value = 20001
$ws.Range("K19").Value = 2450
$ws.Range("L19").Value = 20001
$ws.Range("M19").Value = -2162
$ws.Range("N19").Value = -20577
$ws.Range("H22").Value = 7448.316
$ws.Range("I22").Value = 4875
$ws.Range("J22").Value = 11859.714
$ws.Range("K22").Value = 4875
$ws.Range("L22").Value = 11859.714
$ws.Range("M22").Value = -4346
$ws.Range("N22").Value = -12917.714
$ws.Range("H23").Value = 7671.3335
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 7671.3335
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 7671.3335
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -8117.3335
$ws.Range("H24").Value = 5015500
$ws.Range("J24").Value = 31000
$ws.Range("L24").Value = 31000
$ws.Range("N24").Value = -31346
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("H28").Value = 20015
$ws.Range("J28").Value = 20015
$ws.Range("L28").Value = 20015
$ws.Range("N28").Value = -20399

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2099.6155
$ws.Range("I22").Value = 1024.4166
$ws.Range("J22").Value = 15002
$ws.Range("K22").Value = 1024.4166
$ws.Range("L22").Value = 15002
$ws.Range("M22").Value = -729.4166
$ws.Range("N22").Value = -15592
$ws.Range("H27").Value = 2099.6155
$ws.Range("I27").Value = 1024.4166
$ws.Range("J27").Value = 15002
$ws.Range("K27").Value = 1024.4166
$ws.Range("L27").Value = 15002
$ws.Range("M27").Value = -917.4166
$ws.Range("N27").Value = -15216
$ws.Range("H45").Value = 5900
$ws.Range("I45").Value = 5900
$ws.Range("K45").Value = 5900
$ws.Range("M45").Value = -5493
$ws.Range("H46").Value = 1781.1177
$ws.Range("J46").Value = 1823.6875
$ws.Range("L46").Value = 1823.6875
$ws.Range("N46").Value = -2199.6875

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H103").Value = 20000
$ws.Range("J103").Value = 20000
$ws.Range("L103").Value = 20000
$ws.Range("N103").Value = -22344
$ws.Range("H113").Value = 2798.647
$ws.Range("I113").Value = 867.8570999999999
$ws.Range("J113").Value = 4150.2
$ws.Range("K113").Value = 2603.5713
$ws.Range("L113").Value = 12450.6
$ws.Range("M113").Value = -433.5712999999996
$ws.Range("N113").Value = -16790.6
